# Updated cryptos list on Sun Mar 19 16:09:11 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# reflects three rank swaps (rows 26/27 and the 36-37-38 rotation) by
# rewriting Coin (B) / Link (C) / Price (D) / Volume(1h) (E) in place
# (the rank index in column A is left untouched, exactly like the diff).
#
# Price values that look like plain decimals (e.g. "0.9993") are entered
# with a leading "'" so Excel keeps them as text (quote-prefixed), matching
# how they were already stored in the workbook (inline/shared strings, not
# numbers). Values that already contain extra punctuation (e.g. the BTC/ETH
# "27.603.80" style prices) are naturally kept as text without the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.603.80"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.798.49"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.68%  "

$ws.Range("D5").Value = "'338.90"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").Value = "'0.9966"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'0.3922"
$ws.Range("E7").Value = "  +3.38%  "

$ws.Range("D8").Value = "'0.3469"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'48.09"
$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").Value = "'1.199"
$ws.Range("E10").Value = "  -1.01%  "

$ws.Range("D11").Value = "'0.07519"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "'0.9964"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "'22.13"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "'6.515"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "1.798.01"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").Value = "'7.170"
$ws.Range("E16").Value = "  +1.32%  "

$ws.Range("D17").Value = "'0.00001103"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").Value = "'0.06706"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "'85.18"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'0.9974"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").Value = "'6.574"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("D23").Value = "27.595.59"
$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("D24").Value = "'12.47"
$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.521"
$ws.Range("E26").Value = "  -2.01%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.33"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").Value = "'1.470"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("D29").Value = "'156.75"
$ws.Range("E29").Value = "  +4.51%  "

$ws.Range("D30").Value = "2.002.17"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'135.61"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").Value = "'4.034"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").Value = "'6.085"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").Value = "'0.08784"

$ws.Range("D35").Value = "'13.16"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02430"
$ws.Range("E36").Value = "  +2.81%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.469"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.609"
$ws.Range("E38").Value = "  -3.41%  "

$ws.Range("D39").Value = "'0.06493"
$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("D40").Value = "'0.6861"
$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D41").Value = "'0.2213"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").Value = "'1.254"
$ws.Range("E42").Value = "  -1.45%  "

$ws.Range("D43").Value = "'8.456"
$ws.Range("E43").Value = "  -4.70%  "

$ws.Range("D44").Value = "'14.74"
$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("D45").Value = "'0.6446"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").Value = "'0.9961"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("D47").Value = "'3.871"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("D48").Value = "'2.149"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("D49").Value = "'132.63"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("D50").Value = "'0.07209"
$ws.Range("E50").Value = "  -0.27%  "

$ws.Range("D51").Value = "'80.33"
$ws.Range("E51").Value = "  +0.57%  "
